# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 386
$ws1.Range("F3").Value = 847
$ws1.Range("F5").Value = 1059
$ws1.Range("F6").Value = 2458
$ws1.Range("F7").Value = 208

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 386
$ws4.Range("F3").Value = 847
$ws4.Range("F7").Value = 1059
$ws4.Range("F8").Value = 2458
$ws4.Range("F10").Value = 208
